# explore china energy data
# Replace "James" / "UoC" on Sheet1 row 3 with "Kimi Ma" / "ZJU",
# and move the active selection to D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Kimi Ma"
$ws.Range("B3").Value = "ZJU"

$ws.Range("D10").Select()
